$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B and C, and of columns D and E,
# for the header row and all data rows (rows 1-4).
for ($r = 1; $r -le 4; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal

    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}

# Swap the column widths to match (B<->C, D<->E).
$bWidth = $ws.Columns.Item(2).ColumnWidth
$cWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $cWidth
$ws.Columns.Item(3).ColumnWidth = $bWidth

$dWidth = $ws.Columns.Item(4).ColumnWidth
$eWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $eWidth
$ws.Columns.Item(5).ColumnWidth = $dWidth

# Update the selection to match the new state (column D selected).
$ws.Range("D1:D1048576").Select()
